$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Despezas")

# Update existing row 2 (income entry)
$ws.Range("A2").Value = "R$ 2000.00"
$ws.Range("B2").Value = "Trabalho"
$ws.Range("C2").Value = "Ganho"

# Add row 3 (second income entry)
$ws.Range("A3").Value = "R$ 1500"
$ws.Range("B3").Value = "Trabalho 2"
$ws.Range("C3").Value = "Ganho"

# Add row 4 (expense entry)
$ws.Range("A4").Value = "R$ 560.66"
$ws.Range("B4").Value = "Mercado"
$ws.Range("C4").Value = "Gasto"

# Add row 5 (expense entry)
$ws.Range("A5").Value = "R$ 145.60"
$ws.Range("B5").Value = "Carro"
$ws.Range("C5").Value = "Gasto"
